$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -ne $null -and $val.StartsWith("System, ") -and -not $val.Contains("backdoor")) {
        $rest = $val.Substring(8)
        $newVal = $rest + ", System"
        $cell.Value = $newVal
    }
}
